$d = $word.ActiveDocument

# 1. Update the date in the title: July 2 -> July 3, 2022
$d.Content.Find.Execute(
    "Plant Selection Planner, July 2, 2022", $false, $false, $false, $false, $false,
    $true, 1, $false, "Plant Selection Planner, July 3, 2022", 2) | Out-Null

# 2. Expand the tech-stack sentence
$d.Content.Find.Execute(
    "Foundation, Node.js, Express.js, Handlebars.js, MySQL, Sequelize, nodemailer",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Foundation for sites, Node.js, Express.js, Handlebars.js, MySQL, Sequelize, and nodemailer",
    2) | Out-Null

# 3. Remove the "Submit the application as project 2." bullet paragraph entirely
$paras = $d.Paragraphs
$count = $paras.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*Submit the application as project 2.*") {
        $p.Range.Delete()
        break
    }
}

# 4. Extend the "Front" wireframe description
$d.Content.Find.Execute(
    "Front: The front page contains a brief description of the application and buttons to choose whether to log in, sign up as a gardener or sign up as a nursery manager.",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Front: The front page contains a brief description of the application and buttons to choose whether to log in, sign up as a gardener or sign up as a nursery manager.  The user can search the plant database but cannot save any results without logging in or signing up.",
    2) | Out-Null

# 5. Trim the "as_nursery_manager" wireframe description
$d.Content.Find.Execute(
    "as_nursery_manager: If the user chooses to sign up as a nursery manager he gets a page which lets him enter his desired username, password and email address.  In the future we will also let the nursery manager choose which nursery he is managing, and have a way to create new nurseries.",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "as_nursery_manager: If the user chooses to sign up as a nursery manager he gets a page which lets him enter his desired username, password and email address. ",
    2) | Out-Null

Write-Output "done"
